$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = "200c320a25b192e3c83440f334527e01"
$ws.Range("B89").Value = "540c6e9b1efc86a7027d6bfbd80c73c2"
$ws.Range("B99").Value = "3ed806b97270274a88c3d0a88769021f"
$ws.Range("B110").Value = "1cbee20c6dd597308e23e402c1cb3429"
$ws.Range("B126").Value = "0e7449a6be04ef7efd69afaf0df094cc"
$ws.Range("B154").Value = "7883f0f152cc9d9bb5a1fc710f211227"
$ws.Range("B160").Value = "86c3466b53645a70143a60d23010a457"
$ws.Range("B278").Value = "ff0cdaad1bb498b10fd0b974320bdfa6"
$ws.Range("B335").Value = "ce0d246ac8e46bde9469712017fd6d68"
$ws.Range("B420").Value = "bf3569543f5afe0bd329968445d710df"
$ws.Range("B542").Value = "b526e2e952a95b9a09ec2a8738f95769"
$ws.Range("B561").Value = "5cbb749084cfb11e073fabbd9fa5cca4"
$ws.Range("B580").Value = "a7bcf87a3faf7a525f8737330e459fae"
$ws.Range("B592").Value = "2a0370be441331729a17ae4b1bdd77b2"
$ws.Range("B688").Value = "02796346b86ff6d9d6c7fce4bac0cac5"
$ws.Range("B693").Value = "2d3d3d86d21bacb7bbb70fb06d396780"
$ws.Range("B711").Value = "2bbbc64dc8be0d94d0befb3fe111fabd"
$ws.Range("B776").Value = "ec7cbf44da2741d451e3a0d8eb8e7bff"
$ws.Range("B819").Value = "19e459ae140fd3ca9c68c0372a062362"
$ws.Range("B823").Value = "ce02acf55c77ea096712c1a555e3035c"
$ws.Range("B824").Value = "31a7eec50e7a0a340aa5949d03d55669"
$ws.Range("B833").Value = "138c1287037ebf103f817fe612d3f27d"
$ws.Range("B835").Value = "820a409f29375b7c62388a0b687f0f64"
$ws.Range("B870").Value = "2868f8250a17e53d0e7b5226a008fd5f"
